$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 144, shifting existing rows
# 144-178 down to 145-179 (mirrors the diff's net effect of adding one
# new price record for "Poroto verde" right after row 143).
$ws.Rows(144).Insert()

# Populate the newly inserted row 144 with the new record's data.
$ws.Range("A144").Value = 5
$ws.Range("B144").Value = "Macroferia Regional de Talca"
$ws.Range("C144").Value = "Maule"
$ws.Range("D144").Value = 44855
$ws.Range("E144").Value = 7
$ws.Range("F144").Value = 100112031
$ws.Range("G144").Value = "Poroto verde"
$ws.Range("H144").Value = "Sin especificar"
$ws.Range("I144").Value = "Primera"
$ws.Range("J144").Value = 150
$ws.Range("K144").Value = 28000
$ws.Range("L144").Value = 28000
$ws.Range("M144").Value = 28000
$ws.Range("N144").Value = "`$/malla 25 kilos"
$ws.Range("O144").Value = "Perú"
$ws.Range("P144").Value = 1120
$ws.Range("Q144").Value = 25
$ws.Range("R144").Value = "Hortaliza"
